$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.835.02"
$ws.Range("E2").Value = "  +0.26%  "

$ws.Range("D3").Value = "1.642.36"
$ws.Range("E3").Value = "  +0.07%  "

$ws.Range("E4").Value = "  -0.65%  "

$ws.Range("D5").Value = "'216.21"

$ws.Range("E6").Value = "  +1.14%  "

$ws.Range("E7").Value = "  -0.74%  "

$ws.Range("E8").Value = "  +1.24%  "

$ws.Range("E9").Value = "  -0.24%  "

$ws.Range("D10").Value = "'19.80"
$ws.Range("E10").Value = "  +4.09%  "

$ws.Range("D11").Value = "'0.0844"
$ws.Range("E11").Value = "  -0.19%  "

$ws.Range("D12").Value = "1.872.83"
$ws.Range("E12").Value = "  +0.07%  "

$ws.Range("D13").Value = "1.652.45"
$ws.Range("E13").Value = "  +0.90%  "

$ws.Range("E14").Value = "  -0.17%  "

$ws.Range("E15").Value = "  +0.99%  "

$ws.Range("D16").Value = "'66.37"
$ws.Range("E16").Value = "  +2.95%  "

$ws.Range("D17").Value = "26.849.16"
$ws.Range("E17").Value = "  +0.33%  "

$ws.Range("E18").Value = "  +0.89%  "

$ws.Range("D19").Value = "'217.38"
$ws.Range("E19").Value = "  +3.24%  "

$ws.Range("E20").Value = "  -0.72%  "

$ws.Range("E21").Value = "  +1.28%  "

$ws.Range("E22").Value = "  +7.06%  "

$ws.Range("D23").Value = "'2.43"
$ws.Range("E23").Value = "  +5.16%  "

$ws.Range("D24").Value = "'9.15"
$ws.Range("E24").Value = "  -0.78%  "

$ws.Range("D25").Value = "'145.63"
$ws.Range("E25").Value = "  -1.13%  "

$ws.Range("E26").Value = "  -0.87%  "

$ws.Range("E27").Value = "  +4.38%  "

$ws.Range("E28").Value = "  +0.82%  "

$ws.Range("E29").Value = "  +1.97%  "

$ws.Range("D30").Value = "'0.0511"
$ws.Range("E30").Value = "  +2.19%  "

$ws.Range("E31").Value = "  -0.34%  "

$ws.Range("D32").Value = "'3.36"
$ws.Range("E32").Value = "  +0.76%  "

$ws.Range("E33").Value = "  +1.19%  "

$ws.Range("E34").Value = "  +1.85%  "

$ws.Range("E35").Value = "  -0.45%  "

$ws.Range("D36").Value = "1.241.94"
$ws.Range("E36").Value = "  -2.29%  "

$ws.Range("E37").Value = "  -0.04%  "

$ws.Range("E38").Value = "  +2.82%  "

$ws.Range("E39").Value = "  +3.62%  "

$ws.Range("E40").Value = "  -0.72%  "

$ws.Range("E41").Value = "  +0.19%  "

$ws.Range("E42").Value = "  +2.08%  "

$ws.Range("D43").Value = "1.785.59"
$ws.Range("E43").Value = "  +0.22%  "

$ws.Range("D45").Value = "'60.83"
$ws.Range("E45").Value = "  +1.24%  "

$ws.Range("D46").Value = "'91.52"
$ws.Range("E46").Value = "  +0.34%  "

$ws.Range("E47").Value = "  +0.26%  "

$ws.Range("E48").Value = "  +1.28%  "

$ws.Range("E49").Value = "  -1.02%  "

$ws.Range("E50").Value = "  +1.50%  "

$ws.Range("D51").Value = "'7.55"
$ws.Range("E51").Value = "  +0.88%  "
